# LabsGrilleDeCorrection.xlsx - grading pass ("correction - Mathis fait")
# Fills in the "note" (C) column with the points actually awarded for each
# requirement, and the "Commentaires" (D) column with feedback notes for a
# few rows. Totals in C16/C17 are formulas and recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Points obtained per requirement (column C), mirrored from the max points
# in column B except where a deduction comment was left (rows 2 and 8).
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 8
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 10
$ws.Range("C6").Value = 10
$ws.Range("C7").Value = 10
$ws.Range("C8").Value = 8
$ws.Range("C9").Value = 10
$ws.Range("C10").Value = 5
$ws.Range("C11").Value = 10
$ws.Range("C12").Value = 10
$ws.Range("C13").Value = 10
$ws.Range("C14").Value = 10
$ws.Range("C15").Value = 5

# Comments left in column D for specific requirements.
$ws.Range("D2").Value = "Attention: logo de page de produit est déformé"
$ws.Range("D5").Value = "préférable de centrer les produits, au lieu de laisser 3/4 de la page vide"
$ws.Range("D8").Value = "manque une section de jobs pour George"
$ws.Range("D13").Value = "bon commentaires bien utile en html, mais manque de commentaires en css"

# Apply the same wrap-text style used throughout column A/D to every cell in
# D2:D15 (including the ones left blank) so formatting stays consistent.
$ws.Range("D2:D15").WrapText = $true

# Restore the selection to the cell the author ended up on.
$ws.Range("D16").Select()
